# Daily attendance processing - swap the order of the "Recorded By" names
# in column G from "System, <email>" to "<email>, System" for every row
# where that exact pattern still appears.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

for ($i = 1; $i -le $lastRow; $i++) {
    $cell = $ws.Cells.Item($i, 7)
    if ($cell.Text -eq $oldValue) {
        $cell.Value = $newValue
    }
}
